$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.26%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.81%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.915"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.71%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07322"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.46%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.228"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'22.27%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.739"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.54%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.721"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.26%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9040"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.75%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09101"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'18.23%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1694"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.42%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08256"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.29%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03125"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.88%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09928"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.82%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001503"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.57%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005709"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.15%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.535"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.90%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.078"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.11%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3328"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.31%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-2.43%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.192"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.48%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2102"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-12.01%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04517"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.06%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.44%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004154"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'3.60%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D27").Value = "'0.0003397"
$ws.Range("D27").Style = "Normal"
$ws.Range("D39").Value = "'0.01572"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.67%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04445"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.04%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007336"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.60%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009555"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-4.81%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1326"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.56%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002292"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'11.17%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008322"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-9.48%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006118"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.73%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.413"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'7.45%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-33.27%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.01%"
$ws.Range("E51").Style = "Normal"
